$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Prueba "
$ws.Range("B12").Value = 4242

# "05/02/2026" looks like a date, and a plain .Value assignment would get
# auto-parsed into a date serial (changing both the stored type and the
# cell's number format/style). Route it through a throwaway formula cell
# so it lands as literal text, then paste-special (values only) into the
# target cell so the destination keeps its own existing style untouched.
$ws.Range("G12").Formula = '="05/02/2026"'
$ws.Range("G12").Copy()
$ws.Range("C12").PasteSpecial(-4163)
$ws.Range("G12").Clear()

$ws.Range("D12").Value = 54524224
$ws.Range("E12").Value = "nike"
$ws.Range("F12").Value = "45919158"

$ws.Range("G12").Select()
